$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28").Value = "AutorizationRequired anotacion en metodos privados, para poder validar antes"
$ws.Range("B28").Value = "Lucas"

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A29").Select()
